# Applies updated TPM-derived values to Cd28-Cd86 LR-pairs sheet (rows 2-16)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 6.913788666666666
$ws.Range("H2").Value = 20.741366
$ws.Range("I2").Value = 0.4464851245108818
$ws.Range("J2").Value = 0.4464851245108818
$ws.Range("O2").Value = 0.0004487986043333653
$ws.Range("P2").Value = 0.0004487986043333653
$ws.Range("Q2").Value = 0.643127535562
$ws.Range("R2").Value = 5.788147820058
$ws.Range("S2").Value = 0.0002003819007360926
$ws.Range("T2").Value = 0.0002003819007360926

# Row 3
$ws.Range("G3").Value = 6.913788666666666
$ws.Range("H3").Value = 20.741366
$ws.Range("I3").Value = 0.4464851245108818
$ws.Range("J3").Value = 0.4464851245108818
$ws.Range("O3").Value = 0.0004208040751208546
$ws.Range("P3").Value = 0.0004208040751208545
$ws.Range("Q3").Value = 0.6030114291217777
$ws.Range("R3").Value = 5.427102862096
$ws.Range("S3").Value = 0.0001878827598750212
$ws.Range("T3").Value = 0.0001878827598750212

# Row 4
$ws.Range("G4").Value = 6.913788666666666
$ws.Range("H4").Value = 20.741366
$ws.Range("I4").Value = 0.4464851245108818
$ws.Range("J4").Value = 0.4464851245108818
$ws.Range("M4").Value = 127.2007906666667
$ws.Range("N4").Value = 381.602372
$ws.Range("O4").Value = 0.6137059085722639
$ws.Range("P4").Value = 0.6137059085722638
$ws.Range("Q4").Value = 879.439384902239
$ws.Range("R4").Value = 7914.954464120152
$ws.Range("S4").Value = 0.2740105590019511
$ws.Range("T4").Value = 0.2740105590019511

# Row 5
$ws.Range("G5").Value = 6.913788666666666
$ws.Range("H5").Value = 20.741366
$ws.Range("I5").Value = 0.4464851245108818
$ws.Range("J5").Value = 0.4464851245108818
$ws.Range("O5").Value = 0.0002679928970479904
$ws.Range("P5").Value = 0.0002679928970479904
$ws.Range("Q5").Value = 0.3840333052786667
$ws.Range("R5").Value = 3.456299747508
$ws.Range("S5").Value = 0.0001196548420065039
$ws.Range("T5").Value = 0.0001196548420065039

# Row 6
$ws.Range("G6").Value = 6.913788666666666
$ws.Range("H6").Value = 20.741366
$ws.Range("I6").Value = 0.4464851245108818
$ws.Range("J6").Value = 0.4464851245108818
$ws.Range("M6").Value = 79.830111
$ws.Range("N6").Value = 239.490333
$ws.Range("O6").Value = 0.385156495851234
$ws.Range("P6").Value = 0.385156495851234
$ws.Range("Q6").Value = 551.928516690542
$ws.Range("R6").Value = 4967.356650214878
$ws.Range("S6").Value = 0.1719666460063131
$ws.Range("T6").Value = 0.1719666460063131

# Row 7
$ws.Range("I7").Value = 0.002033181734278123
$ws.Range("J7").Value = 0.002033181734278123
$ws.Range("O7").Value = 0.0004487986043333653
$ws.Range("P7").Value = 0.0004487986043333653
$ws.Range("S7").Value = 0.0000009124891247001131
$ws.Range("T7").Value = 0.000000912489124700113

# Row 8
$ws.Range("I8").Value = 0.002033181734278123
$ws.Range("J8").Value = 0.002033181734278123
$ws.Range("O8").Value = 0.0004208040751208546
$ws.Range("P8").Value = 0.0004208040751208545
$ws.Range("S8").Value = 0.0000008555711592455208
$ws.Range("T8").Value = 0.0000008555711592455207

# Row 9
$ws.Range("I9").Value = 0.002033181734278123
$ws.Range("J9").Value = 0.002033181734278123
$ws.Range("M9").Value = 127.2007906666667
$ws.Range("N9").Value = 381.602372
$ws.Range("O9").Value = 0.6137059085722639
$ws.Range("P9").Value = 0.6137059085722638
$ws.Range("Q9").Value = 4.004747293085778
$ws.Range("R9").Value = 36.04272563777199
$ws.Range("S9").Value = 0.001247775643527687
$ws.Range("T9").Value = 0.001247775643527687

# Row 10
$ws.Range("I10").Value = 0.002033181734278123
$ws.Range("J10").Value = 0.002033181734278123
$ws.Range("O10").Value = 0.0002679928970479904
$ws.Range("P10").Value = 0.0002679928970479904
$ws.Range("S10").Value = 0.0000005448782631942515
$ws.Range("T10").Value = 0.0000005448782631942515

# Row 11
$ws.Range("I11").Value = 0.002033181734278123
$ws.Range("J11").Value = 0.002033181734278123
$ws.Range("M11").Value = 79.830111
$ws.Range("N11").Value = 239.490333
$ws.Range("O11").Value = 0.385156495851234
$ws.Range("P11").Value = 0.385156495851234
$ws.Range("Q11").Value = 2.513344604687
$ws.Range("R11").Value = 22.620101442183
$ws.Range("S11").Value = 0.0007830931522032966
$ws.Range("T11").Value = 0.0007830931522032966

# Row 12
$ws.Range("G12").Value = 8.539652666666667
$ws.Range("H12").Value = 25.618958
$ws.Range("I12").Value = 0.55148169375484
$ws.Range("J12").Value = 0.55148169375484
$ws.Range("O12").Value = 0.0004487986043333653
$ws.Range("P12").Value = 0.0004487986043333653
$ws.Range("Q12").Value = 0.7943670307060001
$ws.Range("R12").Value = 7.149303276354
$ws.Range("S12").Value = 0.0002475042144725726
$ws.Range("T12").Value = 0.0002475042144725726

# Row 13
$ws.Range("G13").Value = 8.539652666666667
$ws.Range("H13").Value = 25.618958
$ws.Range("I13").Value = 0.55148169375484
$ws.Range("J13").Value = 0.55148169375484
$ws.Range("O13").Value = 0.0004208040751208546
$ws.Range("P13").Value = 0.0004208040751208545
$ws.Range("Q13").Value = 0.7448171193831111
$ws.Range("R13").Value = 6.703354074448
$ws.Range("S13").Value = 0.0002320657440865878
$ws.Range("T13").Value = 0.0002320657440865878

# Row 14
$ws.Range("G14").Value = 8.539652666666667
$ws.Range("H14").Value = 25.618958
$ws.Range("I14").Value = 0.55148169375484
$ws.Range("J14").Value = 0.55148169375484
$ws.Range("M14").Value = 127.2007906666667
$ws.Range("N14").Value = 381.602372
$ws.Range("O14").Value = 0.6137059085722639
$ws.Range("P14").Value = 0.6137059085722638
$ws.Range("Q14").Value = 1086.250571218709
$ws.Range("R14").Value = 9776.255140968377
$ws.Range("S14").Value = 0.3384475739267851
$ws.Range("T14").Value = 0.338447573926785

# Row 15
$ws.Range("G15").Value = 8.539652666666667
$ws.Range("H15").Value = 25.618958
$ws.Range("I15").Value = 0.55148169375484
$ws.Range("J15").Value = 0.55148169375484
$ws.Range("O15").Value = 0.0002679928970479904
$ws.Range("P15").Value = 0.0002679928970479904
$ws.Range("Q15").Value = 0.4743435470226667
$ws.Range("R15").Value = 4.269091923204
$ws.Range("S15").Value = 0.0001477931767782922
$ws.Range("T15").Value = 0.0001477931767782922

# Row 16
$ws.Range("G16").Value = 8.539652666666667
$ws.Range("H16").Value = 25.618958
$ws.Range("I16").Value = 0.55148169375484
$ws.Range("J16").Value = 0.55148169375484
$ws.Range("M16").Value = 79.830111
$ws.Range("N16").Value = 239.490333
$ws.Range("O16").Value = 0.385156495851234
$ws.Range("P16").Value = 0.385156495851234
$ws.Range("Q16").Value = 681.7214202814461
$ws.Range("R16").Value = 6135.492782533014
$ws.Range("S16").Value = 0.2124067566927175
$ws.Range("T16").Value = 0.2124067566927175
